$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply a thick red box border to every cell in the used range (A1:G13),
# matching the "red outline" formatting added to the table in the update.
$rng = $ws.Range("A1:G13")
$rng.Borders.Weight = 4
$rng.Borders.Color = 255
